# Apply updated cryptocurrency price/volume data (commit: "Updated cryptos list")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '62.900.48'
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.473.18'
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '414.94'
$ws.Cells.Item(5, 5).Value = '  +1.19%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '130.98'
$ws.Cells.Item(6, 5).Value = '  +1.55%  '
$ws.Cells.Item(7, 5).Value = '  -1.77%  '
$ws.Cells.Item(8, 5).Value = '  +0.05%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.730'
$ws.Cells.Item(9, 5).Value = '  -1.07%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.152'
$ws.Cells.Item(10, 5).Value = '  +7.14%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '42.71'
$ws.Cells.Item(11, 5).Value = '  -2.04%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '9.76'
$ws.Cells.Item(12, 5).Value = '  +3.68%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000223'
$ws.Cells.Item(13, 5).Value = '  -1.87%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '4.027.13'
$ws.Cells.Item(14, 5).Value = '  +1.85%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.141'
$ws.Cells.Item(15, 5).Value = '  -0.27%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '20.56'
$ws.Cells.Item(16, 5).Value = '  -3.83%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '3.487.08'
$ws.Cells.Item(17, 5).Value = '  +1.23%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '12.65'
$ws.Cells.Item(18, 5).Value = '  +0.56%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '1.10'
$ws.Cells.Item(19, 5).Value = '  +0.93%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '62.806.48'
$ws.Cells.Item(20, 5).Value = '  +1.39%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '473.88'
$ws.Cells.Item(21, 5).Value = '  -0.13%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '90.74'
$ws.Cells.Item(22, 5).Value = '  -2.34%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '3.31'
$ws.Cells.Item(23, 5).Value = '  +2.96%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '13.22'
$ws.Cells.Item(24, 5).Value = '  -0.08%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '10.56'
$ws.Cells.Item(25, 5).Value = '  +13.17%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.33'
$ws.Cells.Item(26, 5).Value = '  +0.33%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '33.51'
$ws.Cells.Item(27, 5).Value = '  -0.08%  '
$ws.Cells.Item(28, 5).Value = '  +0.57%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.59'
$ws.Cells.Item(29, 5).Value = '  -0.73%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '12.18'
$ws.Cells.Item(30, 5).Value = '  +0.89%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.69'
$ws.Cells.Item(31, 5).Value = '  -1.06%  '
$ws.Cells.Item(32, 5).Value = '  -0.91%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.113'
$ws.Cells.Item(33, 5).Value = '  -1.42%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '41.05'
$ws.Cells.Item(34, 5).Value = '  -3.18%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.05%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '58.50'
$ws.Cells.Item(36, 5).Value = '  +8.48%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.0489'
$ws.Cells.Item(37, 5).Value = '  -3.13%  '
$ws.Cells.Item(38, 5).Value = '  +0.03%  '
$ws.Cells.Item(39, 2).Value = 'WEMIXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.82'
$ws.Cells.Item(39, 5).Value = '  +8.50%  '
$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.06'
$ws.Cells.Item(40, 5).Value = '  +3.40%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '150.00'
$ws.Cells.Item(41, 5).Value = '  +3.83%  '
$ws.Cells.Item(42, 5).Value = '  -1.34%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.322'
$ws.Cells.Item(43, 5).Value = '  +0.55%  '
$ws.Cells.Item(44, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.34'
$ws.Cells.Item(44, 5).Value = '  -2.21%  '
$ws.Cells.Item(45, 2).Value = 'NEARProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '4.44'
$ws.Cells.Item(45, 5).Value = '  +0.75%  '
$ws.Cells.Item(46, 5).Value = '  +2.30%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0₃0577'
$ws.Cells.Item(47, 5).Value = '  +30.03%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.39'
$ws.Cells.Item(48, 5).Value = '  +10.44%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '16.47'
$ws.Cells.Item(49, 5).Value = '  -1.50%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '22.21'
$ws.Cells.Item(50, 5).Value = '  -1.78%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.143'
$ws.Cells.Item(51, 5).Value = '  -4.76%  '
